# Apply the change described by the diff:
# - Update cell C4 value from 13 to 1.3 (pie_threshold_range Max)
# - Update the active selection on the sheet to C4
# - Resize the saved workbook window to windowWidth=10170 / windowHeight=5340

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data value in C4
$ws.Range("C4").Value = 1.3

# Update the active cell / selection to C4
$ws.Range("C4").Select()

# Reflect the new (smaller) window size recorded in the workbook view
$excel.ActiveWindow.Width = 10170
$excel.ActiveWindow.Height = 5340
